# "Update tabla de mortalidad"
# The values in column B (rows 2-63) had their decimal point removed
# (e.g. 74.245 -> 74245, 62.21 -> 6221, 53.75 -> 5375), turning the
# mortality-rate figures into plain integers. This mirrors exactly what
# happened in the original workbook edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 74245
    3  = 71058
    4  = 6221
    5  = 62212
    6  = 61339
    7  = 60087
    8  = 58712
    9  = 57588
    10 = 56782
    11 = 56169
    12 = 55947
    13 = 55177
    14 = 54695
    15 = 54112
    16 = 5375
    17 = 53274
    18 = 53207
    19 = 53936
    20 = 53266
    21 = 53511
    22 = 52789
    23 = 51807
    24 = 50417
    25 = 49116
    26 = 47159
    27 = 45676
    28 = 44256
    29 = 42979
    30 = 42768
    31 = 42701
    32 = 42873
    33 = 43263
    34 = 43675
    35 = 44054
    36 = 44536
    37 = 45179
    38 = 45913
    39 = 46664
    40 = 47285
    41 = 47825
    42 = 48497
    43 = 4924
    44 = 51572
    45 = 50102
    46 = 50196
    47 = 50343
    48 = 50544
    49 = 50821
    50 = 51025
    51 = 51066
    52 = 51081
    53 = 51437
    54 = 5217
    55 = 52941
    56 = 53687
    57 = 543
    58 = 58565
    59 = 54998
    60 = 55704
    61 = 55515
    62 = 55659
    63 = 55717
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 2).Value = $newValues[$row]
}

# Match the final selection recorded in the workbook after the edit.
$ws.Range("B2:B63").Select()
